$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing "Run 50" column (BA): this also removes the old
# "Mean" header that lived in BA1, and shifts nothing else.
$ws.Columns("BA").Delete()

# Column A header: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Column A body: generation counts -> normalized MaxFES fractions
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Column AZ (now the last column, formerly "Run 50") becomes the "Mean"
# column again, with freshly recomputed values.
$ws.Range("AZ1").Value = "Mean"
$ws.Range("AZ2").Value = 13.32279677
$ws.Range("AZ3").Value = 12.80457169
$ws.Range("AZ4").Value = 11.36591244
$ws.Range("AZ5").Value = 9.57857381
$ws.Range("AZ6").Value = 8.71202011
$ws.Range("AZ7").Value = 8.23425678
$ws.Range("AZ8").Value = 7.44151508
$ws.Range("AZ9").Value = 6.6492936
$ws.Range("AZ10").Value = 5.93668516
$ws.Range("AZ11").Value = 5.1160913
$ws.Range("AZ12").Value = 4.38421681
$ws.Range("AZ13").Value = 3.58882173
$ws.Range("AZ14").Value = 2.8582857
